$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# "Updated XLSX to 0.6.1 / Removed extra lines"
# Drop the three retiring config rows (PowerBudget / DHCPSnooping /
# Temperature): clear the Variable (A) and Description (C) cells
# entirely, and blank the Value (B) cells while leaving their existing
# formatting (style) in place - matches <c r="B20" s="7"/> etc.
$ws.Range("A20:C22").ClearContents() | Out-Null

# Restore the recorded view/selection state for the sheet: scrolled so
# row 7 is the top-left visible cell, with A20:C23 selected.
$ws.Range("A7").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A20:C23").Select() | Out-Null
